# Adds a new "Effective Communicator" glossary section to the "Key Terms"
# sheet: a standalone label in D2, plus six new term/definition rows
# (7-12) formatted like the existing alternating-style rows, a widened
# column B, and the final active selection at B14 — matching the
# author's "Add files via upload" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New standalone label cell ---
$ws.Range("D2").Value = "Effective Communicator"

# --- New glossary rows (terms in column A, definitions in column B) ---
$ws.Range("A7").Value  = "Tone"
$ws.Range("B7").Value  = "The style, emotion, or attitude conveyed in the language. Setting the tone guides the AI's response"

$ws.Range("A8").Value  = "Persona"
$ws.Range("B8").Value  = "A specific identity or perspective the AI should adopt when generating the text. Providing a persona affects the tone and language"

$ws.Range("A9").Value  = "Refine"
$ws.Range("B9").Value  = "Iteratively improve a prompt and response, building on previous interactions rather than starting over each time"

$ws.Range("A10").Value = "Feedback"
$ws.Range("B10").Value = "Asking the AI for suggestions on how to provide better prompts and context. Incorporating this into subsequent prompts"

$ws.Range("A11").Value = "Limitations"
$ws.Range("B11").Value = "Constraints on the AI's knowledge and capabilities. It cannot access real-time info or personal details"

$ws.Range("A12").Value = "Instructions"
$ws.Range("B12").Value = "Clear, direct, and specific requests within a promt. Vague instructions provide vague responses"

# --- Formatting: mirror the existing table's alternating row styles ---
# Column A odd rows (1,3,5,...) use the "yellow header" style; copy it
# from A1 onto the new odd rows 7, 9, 11.
$ws.Range("A1").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A11").PasteSpecial(-4122)

# Column A even rows (2,4,6,...) use the "green header" style; copy it
# from A2 onto the new even rows 8, 10, 12.
$ws.Range("A2").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A12").PasteSpecial(-4122)

# Column B uses the plain bordered style throughout; copy it from B1
# onto the new rows 7-12 in one shot.
$ws.Range("B1").Copy()
$ws.Range("B7:B12").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Column B widens to fit the new, longer definitions ---
$ws.Columns.Item(2).ColumnWidth = 109.346

# --- Final selection left on B14, just past the new table ---
$null = $ws.Range("B14").Select()
